# Updates cryptos list data (prices & 1h volume %) per the
# "Updated cryptos list on Wed May 17 08:33:09 UTC 2023 with GitHub Actions" commit.
#
# Price (column D) and Volume(1h) (column E) cells are stored as literal text in the
# source workbook (many "prices" are not valid Excel numbers, e.g. "26.838.89", and the
# percentages carry padding spaces). Plain `.Value = "..."` assignment lets Excel
# auto-coerce number-looking strings (e.g. "1.001") into real numbers, which would change
# the stored cell type/value. Set-TextValue forces the text interpretation by flipping the
# cell to the Text number format for the write, then restores the original style so no
# visible formatting changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2
Set-TextValue $ws.Range("D2") "26.820.71"
Set-TextValue $ws.Range("E2") "  -1.94%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.808.53"
Set-TextValue $ws.Range("E3") "  -1.22%  "

# Row 4
Set-TextValue $ws.Range("E4") "  -0.42%  "

# Row 5
Set-TextValue $ws.Range("D5") "309.84"
Set-TextValue $ws.Range("E5") "  -1.48%  "

# Row 6
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  -0.33%  "

# Row 7
Set-TextValue $ws.Range("D7") "0.4476"
Set-TextValue $ws.Range("E7") "  +5.27%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3660"
Set-TextValue $ws.Range("E8") "  -1.10%  "

# Row 9
Set-TextValue $ws.Range("D9") "0.07257"
Set-TextValue $ws.Range("E9") "  +0.00%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.8505"
Set-TextValue $ws.Range("E10") "  -1.90%  "

# Row 11
Set-TextValue $ws.Range("D11") "20.62"
Set-TextValue $ws.Range("E11") "  -2.46%  "

# Row 12
Set-TextValue $ws.Range("D12") "1.808.92"
Set-TextValue $ws.Range("E12") "  -1.21%  "

# Row 13
Set-TextValue $ws.Range("D13") "6.593"
Set-TextValue $ws.Range("E13") "  -2.23%  "

# Row 14
Set-TextValue $ws.Range("D14") "0.07074"
Set-TextValue $ws.Range("E14") "  -0.40%  "

# Row 15
Set-TextValue $ws.Range("D15") "5.287"
Set-TextValue $ws.Range("E15") "  -0.56%  "

# Row 16
Set-TextValue $ws.Range("D16") "91.05"
Set-TextValue $ws.Range("E16") "  +1.92%  "

# Row 17
Set-TextValue $ws.Range("D17") "1.002"
Set-TextValue $ws.Range("E17") "  -0.43%  "

# Row 18
Set-TextValue $ws.Range("D18") "0.000008699"
Set-TextValue $ws.Range("E18") "  -1.95%  "

# Row 19
Set-TextValue $ws.Range("D19") "1.001"
Set-TextValue $ws.Range("E19") "  -0.29%  "

# Row 20
Set-TextValue $ws.Range("D20") "14.82"
Set-TextValue $ws.Range("E20") "  -1.92%  "

# Row 21
Set-TextValue $ws.Range("D21") "26.846.72"
Set-TextValue $ws.Range("E21") "  -1.98%  "

# Row 22
Set-TextValue $ws.Range("D22") "5.132"
Set-TextValue $ws.Range("E22") "  -0.12%  "

# Row 23
Set-TextValue $ws.Range("E23") "  -0.87%  "

# Row 24
Set-TextValue $ws.Range("D24") "1.984"
Set-TextValue $ws.Range("E24") "  -0.90%  "

# Row 25
Set-TextValue $ws.Range("D25") "151.20"
Set-TextValue $ws.Range("E25") "  -1.11%  "

# Row 26
Set-TextValue $ws.Range("D26") "2.217"
Set-TextValue $ws.Range("E26") "  +0.89%  "

# Row 27
Set-TextValue $ws.Range("D27") "18.36"
Set-TextValue $ws.Range("E27") "  -0.55%  "

# Row 28
Set-TextValue $ws.Range("D28") "5.185"
Set-TextValue $ws.Range("E28") "  -1.37%  "

# Row 29
Set-TextValue $ws.Range("D29") "115.89"
Set-TextValue $ws.Range("E29") "  -0.63%  "

# Row 30
Set-TextValue $ws.Range("D30") "0.08811"
Set-TextValue $ws.Range("E30") "  -0.85%  "

# Row 31
Set-TextValue $ws.Range("D31") "1.169"
Set-TextValue $ws.Range("E31") "  -2.95%  "

# Row 32
Set-TextValue $ws.Range("D32") "0.7449"
Set-TextValue $ws.Range("E32") "  -2.25%  "

# Row 33
Set-TextValue $ws.Range("D33") "2.930"
Set-TextValue $ws.Range("E33") "  +3.78%  "

# Row 34
Set-TextValue $ws.Range("D34") "4.417"
Set-TextValue $ws.Range("E34") "  -1.60%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.001"
Set-TextValue $ws.Range("E35") "  -0.36%  "

# Row 36
Set-TextValue $ws.Range("D36") "1.085"
Set-TextValue $ws.Range("E36") "  -3.57%  "

# Row 37
Set-TextValue $ws.Range("D37") "0.01951"
Set-TextValue $ws.Range("E37") "  -1.65%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.05170"
Set-TextValue $ws.Range("E38") "  -2.38%  "

# Row 39
Set-TextValue $ws.Range("D39") "0.5257"
Set-TextValue $ws.Range("E39") "  +3.21%  "

# Row 40
Set-TextValue $ws.Range("D40") "2.869"
Set-TextValue $ws.Range("E40") "  -0.48%  "

# Row 41
Set-TextValue $ws.Range("D41") "7.065"
Set-TextValue $ws.Range("E41") "  -2.71%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.1684"
Set-TextValue $ws.Range("E42") "  -1.31%  "

# Row 43
Set-TextValue $ws.Range("D43") "0.5160"
Set-TextValue $ws.Range("E43") "  +7.89%  "

# Row 44
Set-TextValue $ws.Range("D44") "8.410"
Set-TextValue $ws.Range("E44") "  -3.54%  "

# Row 45
Set-TextValue $ws.Range("D45") "10.52"
Set-TextValue $ws.Range("E45") "  -1.20%  "

# Row 46
Set-TextValue $ws.Range("D46") "1.948"
Set-TextValue $ws.Range("E46") "  +5.40%  "

# Row 47
Set-TextValue $ws.Range("D47") "105.13"
Set-TextValue $ws.Range("E47") "  -2.67%  "

# Row 48
Set-TextValue $ws.Range("E48") "  -0.39%  "

# Row 49
Set-TextValue $ws.Range("B49") "NEARProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.650"
Set-TextValue $ws.Range("E49") "  -1.37%  "

# Row 50
Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.06319"
Set-TextValue $ws.Range("E50") "  -1.23%  "

# Row 51
Set-TextValue $ws.Range("D51") "0.9112"
Set-TextValue $ws.Range("E51") "  -1.09%  "
